# Refresh the crypto price ("Price", col D) and volume-change ("Volume(1h)",
# col E) figures for every coin row (2-51) per the scheduled GitHub Actions
# data pull.
#
# Column D cells are plain text (some values, e.g. "26.024.17", use dots as
# thousands separators and are not valid numbers; others, e.g. "4.360", are
# numeric-looking but must keep a trailing zero that Excel's normal value
# parser would otherwise strip). To guarantee every D cell lands back as an
# exact, unstyled text literal (matching the original inlineStr cells), we
# write it as a quoted-string formula first and then Copy/PasteSpecial
# (values only, xlPasteValues = -4163) over itself - this freezes the text
# result without leaving the NumberFormat/quote-prefix style footprint that
# `.Value = "..."` or an apostrophe-prefix would otherwise introduce.
#
# Column E cells ("  -3.34%  " style strings) keep their padding spaces, so
# plain `.Value = ` assignment already stays text - no special handling
# needed there.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "=`"26.059.99`""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  -2.61%  "
$ws.Range("D3").Formula = "=`"1.642.94`""
$ws.Range("D3").Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Formula = "=`"215.41`""
$ws.Range("D5").Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").Formula = "=`"0.5079`""
$ws.Range("D6").Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  -2.96%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Formula = "=`"0.2583`""
$ws.Range("D8").Copy()
$ws.Range("D8").PasteSpecial(-4163)
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Formula = "=`"0.06413`""
$ws.Range("D9").Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  -2.36%  "
$ws.Range("E10").Value = "  -3.44%  "
$ws.Range("D11").Formula = "=`"0.07723`""
$ws.Range("D11").Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.13%  "
$ws.Range("D12").Formula = "=`"1.653.66`""
$ws.Range("D12").Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  -2.35%  "
$ws.Range("D13").Formula = "=`"4.261`""
$ws.Range("D13").Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -2.59%  "
$ws.Range("D14").Formula = "=`"1.870.05`""
$ws.Range("D14").Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Formula = "=`"0.5459`""
$ws.Range("D15").Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  -3.54%  "
$ws.Range("D16").Formula = "=`"0.0₅7975`""
$ws.Range("D16").Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Formula = "=`"63.78`""
$ws.Range("D17").Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Formula = "=`"26.068.68`""
$ws.Range("D18").Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  -3.01%  "
$ws.Range("E19").Value = "  +0.04%  "
$ws.Range("D20").Formula = "=`"206.98`""
$ws.Range("D20").Copy()
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = "  -3.95%  "
$ws.Range("D21").Formula = "=`"4.360`""
$ws.Range("D21").Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -4.17%  "
$ws.Range("D22").Formula = "=`"10.01`""
$ws.Range("D22").Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -2.09%  "
$ws.Range("D23").Formula = "=`"5.980`""
$ws.Range("D23").Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +0.60%  "
$ws.Range("D24").Formula = "=`"1.006`""
$ws.Range("D24").Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Formula = "=`"1.969`""
$ws.Range("D25").Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +13.85%  "
$ws.Range("D26").Formula = "=`"142.81`""
$ws.Range("D26").Copy()
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").Formula = "=`"0.1162`""
$ws.Range("D27").Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Formula = "=`"6.878`""
$ws.Range("D28").Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("D29").Formula = "=`"15.83`""
$ws.Range("D29").Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  -0.73%  "
$ws.Range("D30").Formula = "=`"0.05031`""
$ws.Range("D30").Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  -4.86%  "
$ws.Range("D31").Formula = "=`"1.239`""
$ws.Range("D31").Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  -3.07%  "
$ws.Range("D32").Formula = "=`"3.303`""
$ws.Range("D32").Copy()
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = "  -2.72%  "
$ws.Range("E33").Value = "  -2.06%  "
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("D35").Formula = "=`"2.338`""
$ws.Range("D35").Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -2.34%  "
$ws.Range("D36").Formula = "=`"0.9120`""
$ws.Range("D36").Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Formula = "=`"2.656`""
$ws.Range("D37").Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -4.77%  "
$ws.Range("D38").Formula = "=`"0.5696`""
$ws.Range("D38").Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Formula = "=`"1.131.29`""
$ws.Range("D39").Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -3.65%  "
$ws.Range("D40").Formula = "=`"0.01567`""
$ws.Range("D40").Copy()
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = "  -3.38%  "
$ws.Range("D41").Formula = "=`"1.005`""
$ws.Range("D41").Copy()
$ws.Range("D41").PasteSpecial(-4163)
$ws.Range("E41").Value = "  -0.23%  "
$ws.Range("D42").Formula = "=`"2.552`""
$ws.Range("D42").Copy()
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Formula = "=`"5.629`""
$ws.Range("D43").Copy()
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").Formula = "=`"0.8197`""
$ws.Range("D44").Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  -1.65%  "
$ws.Range("D45").Formula = "=`"99.74`""
$ws.Range("D45").Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  -0.22%  "
$ws.Range("D46").Formula = "=`"1.781.33`""
$ws.Range("D46").Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  -2.78%  "
$ws.Range("E47").Value = "  +2.85%  "
$ws.Range("D48").Formula = "=`"0.4527`""
$ws.Range("D48").Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Formula = "=`"1.006`""
$ws.Range("D49").Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -0.09%  "
$ws.Range("D50").Formula = "=`"55.04`""
$ws.Range("D50").Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  -2.78%  "
$ws.Range("D51").Formula = "=`"7.779`""
$ws.Range("D51").Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  -3.08%  "
$excel.CutCopyMode = 0
